# reference_parameters.xlsx column layout fix: columns are reordered to the
# alphabetically-sorted parameter-field order (label, expression, maximum,
# minimum, non_negative, standard_error, value, vary) used by glotaran's
# ParameterDataFrame export, and the "no_defaults" row gets real
# minimum/maximum/non_negative values instead of blank/False placeholders.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -------------------------------------------------------
$ws.Cells.Item(1, 1).Value = "label"
$ws.Cells.Item(1, 2).Value = "expression"
$ws.Cells.Item(1, 3).Value = "maximum"
$ws.Cells.Item(1, 4).Value = "minimum"
$ws.Cells.Item(1, 5).Value = "non_negative"
$ws.Cells.Item(1, 6).Value = "standard_error"
$ws.Cells.Item(1, 7).Value = "value"
$ws.Cells.Item(1, 8).Value = "vary"

# --- Data rows ----------------------------------------------------------
# columns: label, expression, maximum, minimum, non_negative, standard_error, value, vary
$data = @(
    @("pure_list.1",                   "None", "", "", $false, "None", 1,  $true),
    @("pure_list.2",                   "None", "", "", $false, "None", 2,  $true),
    @("list_with_options.1",           "None", "", "", $false, "None", 3,  $false),
    @("list_with_options.2",           "None", "", "", $false, "None", 4,  $false),
    @("verbose_list.all_defaults",     "None", "", "", $false, "None", 5,  $true),
    @("verbose_list.no_defaults",      "None", 1,  -1, $true,  "None", 6,  $false),
    @("verbose_list.expression_only",  '$verbose_list.all_defaults + $verbose_list.no_defaults', "", "", $false, "None", 11, $false)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}
